$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65
$ws.Cells.Item(65, 1).Value = 705567
$ws.Cells.Item(65, 2).Value = 705567
$ws.Cells.Item(65, 3).Value = ""

# Row 66
$ws.Cells.Item(66, 1).Value = 768046
$ws.Cells.Item(66, 2).Value = 768046
$ws.Cells.Item(66, 3).Value = ""
$v1 = @"
1 de Set de 2025
"@
$ws.Cells.Item(66, 4).Value = $v1
$v2 = @"
Naturalização Ordinária
"@
$ws.Cells.Item(66, 5).Value = $v2
$v3 = @"
Indeferimento
"@
$ws.Cells.Item(66, 6).Value = $v3
$v4 = @"
Art. 65, inciso II da Lei nº 13.445/2017; 🚨 REQUERENTE NÃO ESTÁ NO PAÍS - INDEFERIMENTO AUTOMÁTICO; Não anexou item 2
"@
$ws.Cells.Item(66, 7).Value = $v4
$v5 = @"
Indeferimento
"@
$ws.Cells.Item(66, 8).Value = $v5
$v6 = @"
🚨 REQUERENTE NÃO ESTÁ NO PAÍS - INDEFERIMENTO AUTOMÁTICO
"@
$ws.Cells.Item(66, 9).Value = $v6
$v7 = @"
Processo indeferido por não atender aos requisitos
"@
$ws.Cells.Item(66, 10).Value = $v7
$v8 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(66, 11).Value = $v8
$v9 = @"
❌ NÃO ATENDIDO - Prazo de residência não localizado nos campos do sistema
"@
$ws.Cells.Item(66, 12).Value = $v9
$v10 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(66, 13).Value = $v10
$v11 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(66, 14).Value = $v11
$v12 = @"
✅ 75% (3/4)
"@
$ws.Cells.Item(66, 15).Value = $v12
$v13 = @"
6/8
"@
$ws.Cells.Item(66, 16).Value = $v13
$v14 = @"
'75.0%
"@
$ws.Cells.Item(66, 17).Value = $v14
$ws.Cells.Item(66, 17).Style = "Normal"
$v15 = @"
22/11/2025
"@
$ws.Cells.Item(66, 18).Value = $v15
$v16 = @"
14:36:14
"@
$ws.Cells.Item(66, 19).Value = $v16
$v17 = @"
1. Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de Naturalização Ordinária apresento o presente Relatório Opinativo.
2. Conforme registro no Sistema de Tráfego Internacional - STI e no passaporte, o requerente não se encontra em território nacional na data da entrada do processo, impedindo a continuidade do processo.
3. Diante do exposto, encaminhe-se ao Ministério da Justiça e Segurança Pública com opinião DESFAVORÁVEL AO DEFERIMENTO em razão ...
"@
$ws.Cells.Item(66, 20).Value = $v17
$v18 = @"
Não atendeu 2 requisito(s)
"@
$ws.Cells.Item(66, 21).Value = $v18

# Row 67
$ws.Cells.Item(67, 1).Value = 767995
$ws.Cells.Item(67, 2).Value = 767995
$ws.Cells.Item(67, 3).Value = ""
$v19 = @"
1 de Set de 2025
"@
$ws.Cells.Item(67, 4).Value = $v19
$v20 = @"
Naturalização Ordinária
"@
$ws.Cells.Item(67, 5).Value = $v20
$v21 = @"
Indeferimento
"@
$ws.Cells.Item(67, 6).Value = $v21
$v22 = @"
⚠️ AUSÊNCIA DE COLETA BIOMÉTRICA CONSTATADA NO PARECER PF
"@
$ws.Cells.Item(67, 7).Value = $v22
$v23 = @"
Indeferimento
"@
$ws.Cells.Item(67, 8).Value = $v23
$v24 = @"
⚠️ AUSÊNCIA DE COLETA BIOMÉTRICA CONSTATADA NO PARECER PF
"@
$ws.Cells.Item(67, 9).Value = $v24
$v25 = @"
Processo indeferido por não atender aos requisitos
"@
$ws.Cells.Item(67, 10).Value = $v25
$v26 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(67, 11).Value = $v26
$v27 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(67, 12).Value = $v27
$v28 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(67, 13).Value = $v28
$v29 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(67, 14).Value = $v29
$v30 = @"
✅ 100% (4/4)
"@
$ws.Cells.Item(67, 15).Value = $v30
$v31 = @"
8/8
"@
$ws.Cells.Item(67, 16).Value = $v31
$v32 = @"
'100.0%
"@
$ws.Cells.Item(67, 17).Value = $v32
$ws.Cells.Item(67, 17).Style = "Normal"
$v33 = @"
22/11/2025
"@
$ws.Cells.Item(67, 18).Value = $v33
$v34 = @"
14:38:42
"@
$ws.Cells.Item(67, 19).Value = $v34
$v35 = @"
1.	Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de Naturalização Ordinária apresento o presente Relatório Opinativo.
2.	De acordo com a documentação apresentada, o interessado possui capacidade civil. 
3.	A relação de documentos exigidos pela legislação, NÃO foi apresentada integralmente conforme documentos juntados ao processo. NÃO apresentou:
a.	Legalização ou apostilamento do atestado de antecedentes criminais emitido pelo país de origem (o do...
"@
$ws.Cells.Item(67, 20).Value = $v35
$v36 = @"
Não atendeu 1 requisito(s)
"@
$ws.Cells.Item(67, 21).Value = $v36

# Row 68
$ws.Cells.Item(68, 1).Value = 767111
$ws.Cells.Item(68, 2).Value = 767111
$ws.Cells.Item(68, 3).Value = ""
$v37 = @"
29 de Ago de 2025
"@
$ws.Cells.Item(68, 4).Value = $v37
$v38 = @"
Naturalização Ordinária
"@
$ws.Cells.Item(68, 5).Value = $v38
$v39 = @"
Analise Manual
"@
$ws.Cells.Item(68, 6).Value = $v39
$v40 = @"
Art. 65, inciso II da Lei nº 13.445/2017
"@
$ws.Cells.Item(68, 7).Value = $v40
$v41 = @"
Indeferimento
"@
$ws.Cells.Item(68, 8).Value = $v41
$v42 = @"
⚠️ PARECER PF SEM PRAZO DE RESIDÊNCIA ESPECIFICADO
"@
$ws.Cells.Item(68, 9).Value = $v42
$v43 = @"
Processo encaminhado para ANÁLISE MANUAL devido a alerta(s) crítico(s) no parecer da PF ou dados insuficientes para decisão automática.
"@
$ws.Cells.Item(68, 10).Value = $v43
$v44 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(68, 11).Value = $v44
$v45 = @"
❌ NÃO ATENDIDO - Tempo insuficiente: 1.00 anos < 3.95 anos
"@
$ws.Cells.Item(68, 12).Value = $v45
$v46 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(68, 13).Value = $v46
$v47 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(68, 14).Value = $v47
$v48 = @"
✅ 100% (4/4)
"@
$ws.Cells.Item(68, 15).Value = $v48
$v49 = @"
7/8
"@
$ws.Cells.Item(68, 16).Value = $v49
$v50 = @"
'87.5%
"@
$ws.Cells.Item(68, 17).Value = $v50
$ws.Cells.Item(68, 17).Style = "Normal"
$v51 = @"
22/11/2025
"@
$ws.Cells.Item(68, 18).Value = $v51
$v52 = @"
14:40:59
"@
$ws.Cells.Item(68, 19).Value = $v52
$v53 = @"
Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de naturalização, apresento o presente Relatório Opinativo. 
O artigo 65 da Lei 13.445/17 estabelece os seguintes requisitos para o deferimento da naturalização ordinária:
I - ter capacidade civil, segundo a lei brasileira;
II - ter residência em território nacional, pelo prazo mínimo de 4 (quatro) anos;
III - comunicar-se em língua portuguesa, consideradas as condições do naturalizando; e
IV - nã...
"@
$ws.Cells.Item(68, 20).Value = $v53
$v54 = @"
Caso marcado para análise manual (sem decisão automática de deferimento/indeferimento).
"@
$ws.Cells.Item(68, 21).Value = $v54

# Row 69
$ws.Cells.Item(69, 1).Value = 765992
$ws.Cells.Item(69, 2).Value = 765992
$ws.Cells.Item(69, 3).Value = ""
$v55 = @"
27 de Ago de 2025
"@
$ws.Cells.Item(69, 4).Value = $v55
$v56 = @"
Naturalização Ordinária
"@
$ws.Cells.Item(69, 5).Value = $v56
$v57 = @"
Deferimento
"@
$ws.Cells.Item(69, 6).Value = $v57
$ws.Cells.Item(69, 7).Value = ""
$v58 = @"
Indeferimento
"@
$ws.Cells.Item(69, 8).Value = $v58
$v59 = @"
Nenhum
"@
$ws.Cells.Item(69, 9).Value = $v59
$v60 = @"
Processo deferido automaticamente com base na análise de elegibilidade.
"@
$ws.Cells.Item(69, 10).Value = $v60
$v61 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(69, 11).Value = $v61
$v62 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(69, 12).Value = $v62
$v63 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(69, 13).Value = $v63
$v64 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(69, 14).Value = $v64
$v65 = @"
✅ 100% (4/4)
"@
$ws.Cells.Item(69, 15).Value = $v65
$v66 = @"
8/8
"@
$ws.Cells.Item(69, 16).Value = $v66
$v67 = @"
'100.0%
"@
$ws.Cells.Item(69, 17).Value = $v67
$ws.Cells.Item(69, 17).Style = "Normal"
$v68 = @"
22/11/2025
"@
$ws.Cells.Item(69, 18).Value = $v68
$v69 = @"
14:44:28
"@
$ws.Cells.Item(69, 19).Value = $v69
$v70 = @"
RELATóRIO OPINATIVO/PARECER
DATA 08/09/2025
Naturalizando
MOISE STERLING - G2975982: 
1. Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de Naturalização Ordinária apresento o presente Relatório Opinativo/Parecer.
2. A relação de documentos exigidos pela legislação, não foi apresentada integralmente conforme documentos juntados ao processo, visto que o naturalizando não compareceu nem justificou a falta ao agendamento/notificação.
3. Notific...
"@
$ws.Cells.Item(69, 20).Value = $v70
$v71 = @"
Todos os requisitos atendidos segundo a análise automática.
"@
$ws.Cells.Item(69, 21).Value = $v71

# Row 70
$ws.Cells.Item(70, 1).Value = 763365
$ws.Cells.Item(70, 2).Value = 763365
$ws.Cells.Item(70, 3).Value = ""
$v72 = @"
20 de Ago de 2025
"@
$ws.Cells.Item(70, 4).Value = $v72
$v73 = @"
Naturalização Ordinária
"@
$ws.Cells.Item(70, 5).Value = $v73
$v74 = @"
Deferimento
"@
$ws.Cells.Item(70, 6).Value = $v74
$ws.Cells.Item(70, 7).Value = ""
$v75 = @"
Indeferimento
"@
$ws.Cells.Item(70, 8).Value = $v75
$v76 = @"
Nenhum
"@
$ws.Cells.Item(70, 9).Value = $v76
$v77 = @"
Processo deferido automaticamente com base na análise de elegibilidade.
"@
$ws.Cells.Item(70, 10).Value = $v77
$v78 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(70, 11).Value = $v78
$v79 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(70, 12).Value = $v79
$v80 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(70, 13).Value = $v80
$v81 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(70, 14).Value = $v81
$v82 = @"
✅ 100% (4/4)
"@
$ws.Cells.Item(70, 15).Value = $v82
$v83 = @"
8/8
"@
$ws.Cells.Item(70, 16).Value = $v83
$v84 = @"
'100.0%
"@
$ws.Cells.Item(70, 17).Value = $v84
$ws.Cells.Item(70, 17).Style = "Normal"
$v85 = @"
22/11/2025
"@
$ws.Cells.Item(70, 18).Value = $v85
$v86 = @"
14:46:23
"@
$ws.Cells.Item(70, 19).Value = $v86
$v87 = @"
RELATóRIO OPINATIVO/PARECER
DATA 29/08/2025
Naturalizando
JAMESLEY ESTIMABLE - G314924Y: 
1. Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de Naturalização Ordinária apresento o presente Relatório Opinativo/Parecer.
2. A relação de documentos exigidos pela legislação, não foi apresentada integralmente conforme documentos juntados ao processo, visto que o naturalizando não compareceu nem justificou a falta ao agendamento/notificação.
3. Not...
"@
$ws.Cells.Item(70, 20).Value = $v87
$v88 = @"
Todos os requisitos atendidos segundo a análise automática.
"@
$ws.Cells.Item(70, 21).Value = $v88

# Row 71
$v89 = @"
'762641
"@
$ws.Cells.Item(71, 1).Value = $v89
$ws.Cells.Item(71, 1).Style = "Normal"
$v90 = @"
'762641
"@
$ws.Cells.Item(71, 2).Value = $v90
$ws.Cells.Item(71, 2).Style = "Normal"
$v91 = @"
N/A
"@
$ws.Cells.Item(71, 3).Value = $v91
$v92 = @"
18 de Ago de 2025
"@
$ws.Cells.Item(71, 4).Value = $v92
$v93 = @"
Naturalização Ordinária
"@
$ws.Cells.Item(71, 5).Value = $v93
$v94 = @"
Indeferimento
"@
$ws.Cells.Item(71, 6).Value = $v94
$v95 = @"
Art. 65, inciso IV da Lei nº 13.445/2017; ⚠️ AUSÊNCIA DE COLETA BIOMÉTRICA CONSTATADA NO PARECER PF; Não anexou item 4
"@
$ws.Cells.Item(71, 7).Value = $v95
$v96 = @"
Indeferimento
"@
$ws.Cells.Item(71, 8).Value = $v96
$v97 = @"
⚠️ AUSÊNCIA DE COLETA BIOMÉTRICA CONSTATADA NO PARECER PF
"@
$ws.Cells.Item(71, 9).Value = $v97
$v98 = @"
Processo indeferido por não atender aos requisitos
"@
$ws.Cells.Item(71, 10).Value = $v98
$v99 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(71, 11).Value = $v99
$v100 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(71, 12).Value = $v100
$v101 = @"
✅ ATENDIDO
"@
$ws.Cells.Item(71, 13).Value = $v101
$v102 = @"
❌ NÃO ATENDIDO - Antecedentes criminais inválidos ou não anexados
"@
$ws.Cells.Item(71, 14).Value = $v102
$v103 = @"
✅ 75% (3/4)
"@
$ws.Cells.Item(71, 15).Value = $v103
$v104 = @"
6/8
"@
$ws.Cells.Item(71, 16).Value = $v104
$v105 = @"
'75.0%
"@
$ws.Cells.Item(71, 17).Value = $v105
$ws.Cells.Item(71, 17).Style = "Normal"
$v106 = @"
22/11/2025
"@
$ws.Cells.Item(71, 18).Value = $v106
$v107 = @"
14:47:25
"@
$ws.Cells.Item(71, 19).Value = $v107
$v108 = @"
1.	Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de Naturalização Ordinária apresento o presente Relatório Opinativo.
2.	De acordo com a documentação apresentada, o interessado possui capacidade civil. 
3.	A relação de documentos exigidos pela legislação, NÃO foi apresentada integralmente conforme documentos juntados ao processo. NÃO apresentou:
a.	Cópia da Carteira de Registro Nacional Migratório;
b.	Comprovante de situação cadastral do Cadastro ...
"@
$ws.Cells.Item(71, 20).Value = $v108
$v109 = @"
Não atendeu 2 requisito(s)
"@
$ws.Cells.Item(71, 21).Value = $v109

Write-Output "applied edits"
